# Fix category modules and align income model with addIncome flow
#
# The sheet previously held 3 identical "asadullahahmed01@gmail.com" / 233
# expense rows. Replace them with the real income entries ("Drind") and
# drop the now-superfluous 3rd row so the sheet matches the addIncome data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Drind / 7000 / 2026-02-18 06:00:20
$ws.Range("A2").Value = "Drind"
$ws.Range("B2").Value = 7000
$ws.Range("C2").Value = 46071.250231481485

# Row 3: Drind / 3000 / 2026-02-16 06:00:20
$ws.Range("A3").Value = "Drind"
$ws.Range("B3").Value = 3000
$ws.Range("C3").Value = 46069.250231481485

# Remove the now-unused 4th row entirely, shifting the used range/dimension
# back down to A1:C3.
$ws.Rows.Item(4).Delete()
